$d = $word.ActiveDocument

# Locate the "k = zoom in and toggle through vertical image quarters"
# paragraph (new shortcut is inserted right after it, before the
# "i = toggle the user IDs on/off" paragraph).
$kPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "k = zoom in and toggle through vertical image quarters*") {
        $kPara = $p
        break
    }
}

# Duplicate the whole paragraph (incl. paragraph mark) via copy/paste so the
# new paragraph inherits the same paragraph formatting (spacing, fonts) and
# the pasted bold run correctly carries both <w:b/> and <w:bCs/>.
$fullRange = $d.Range($kPara.Range.Start, $kPara.Range.End)
$fullRange.Copy()

$insertPoint = $d.Range($kPara.Range.End, $kPara.Range.End)
$insertPoint.Paste()

# Re-fetch the newly inserted paragraph (right after the original "k" one).
$newPara = $d.Paragraphs.Item($kPara.Index + 1)

# The duplicated paragraph currently reads:
#   "k = zoom in and toggle through vertical image quarters"
# with the leading "k" bold. Turn it into:
#   "a = toggle auto-naming neurons"
# keeping the leading character bold (with its bCs companion).
$boldChar = $d.Range($newPara.Range.Start, $newPara.Range.Start + 1)
$boldChar.Text = "a"

$restStart = $newPara.Range.Start + 1
$restEnd = $newPara.Range.End - 1
$restRange = $d.Range($restStart, $restEnd)
$restRange.Text = " = toggle auto-naming neurons"
